$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.929.81"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Range('E2').Value = '  -0.58%  '
$ws.Range('D3').Value = "'2.215.57"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'257.06"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Range('E5').Value = '  +2.77%  '
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('D7').Value = "'77.71"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Range('E7').Value = '  +3.36%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.08%  '
$ws.Range('D10').Value = "'43.03"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Range('E10').Value = '  +4.63%  '
$ws.Range('E11').Value = '  -2.28%  '
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('D14').Value = "'2.550.98"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Range('E14').Value = '  -1.16%  '
$ws.Range('D15').Value = "'14.45"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Range('E15').Value = '  -1.33%  '
$ws.Range('D16').Value = "'2.220.77"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').Value = "'0.784"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Range('E17').Value = '  -1.02%  '
$ws.Range('D18').Value = "'42.874.01"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').Value = "'0.0000104"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').Value = "'71.09"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').Value = "'5.97"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  +3.34%  '
$ws.Range('D23').Value = "'230.01"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = "'9.32"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Range('E24').Value = '  -4.26%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').Value = "'42.82"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Range('E26').Value = '  +8.89%  '
$ws.Range('D27').Value = "'10.74"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Range('E27').Value = '  -0.60%  '
$ws.Range('E28').Value = '  -2.84%  '
$ws.Range('E29').Value = '  -0.34%  '
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('D31').Value = "'173.92"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').Value = "'20.41"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('D33').Value = "'0.0873"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Range('E33').Value = '  +8.66%  '
$ws.Range('D34').Value = "'5.22"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Range('E34').Value = '  -0.69%  '
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').Value = "'0.0356"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Range('E36').Value = '  +6.96%  '
$ws.Range('E37').Value = '  -2.00%  '
$ws.Range('D38').Value = "'4.38"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('D39').Value = "'13.11"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('D40').Value = "'2.83"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Range('E40').Value = '  +16.72%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').Value = "'0.202"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Range('E42').Value = '  -1.85%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').Value = "'61.39"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Range('E43').Value = '  +3.26%  '
$ws.Range('E44').Value = '  -2.17%  '
$ws.Range('D45').Value = "'0.489"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Range('E45').Value = '  +0.43%  '
$ws.Range('D46').Value = "'103.44"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Range('D47').Value = "'8.45"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').Value = "'0.0969"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Range('E48').Value = '  -2.30%  '
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('D50').Value = "'1.13"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').Value = "'1.47"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Range('E51').Value = '  +21.56%  '
